$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 4 new company/test rows (35-38), following the existing table pattern ---
# Column A and E values are entered in the same interleaved order the source
# workbook's sharedStrings table implies (A35, A36, E35, E36, A37, E37, A38, E38)
# so new shared-string indices line up with the target file.
$ws.Range("A35").Value = "EmpT33"
$ws.Range("A36").Value = "EmpT34"
$ws.Range("E35").Value = "Teste33"
$ws.Range("E36").Value = "Teste34"
$ws.Range("A37").Value = "EmpT35"
$ws.Range("E37").Value = "Teste35"
$ws.Range("A38").Value = "EmpT36"
$ws.Range("E38").Value = "Teste36"

$ws.Range("B35").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 45747

$ws.Range("B36").Value = 4
$ws.Range("C36").Value = 2
$ws.Range("D36").Value = 45747

$ws.Range("B37").Value = 2
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 45747

$ws.Range("B38").Value = 3
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 45747

# Re-apply the same cell formatting used by the row(s) immediately above the
# new block (quote-prefix number style on B/C, date style on D) so the new
# rows look like the rest of the table. Doing this AFTER the values are
# written keeps the values while restoring the format.
$ws.Range("A31:E31").Copy()
$ws.Range("A35:E38").PasteSpecial(-4122)

# --- Update the sheet view: scroll the frozen pane down to the new rows and
#     move the active selection to follow the last-used cell. ---
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E41").Select()

# --- Column F width changed (no data lives in that column; its width was
#     simply re-set by the author, e.g. via a manual resize/auto-fit). ---
$ws.Columns("F").ColumnWidth = 10.833333333333334
